$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("2018")
$ws.Range("C2").Value = 423182.0660165041
$ws.Range("D2").Value = 423639.1464891222
$ws.Range("E2").Value = 763.6294073518379
$ws.Range("F2").Value = 429000
$ws.Range("G2").Value = 430000
$ws.Range("H2").Value = 649
$ws.Range("I2").Value = 1333
$ws.Range("J2").Value = 554.77
$ws.Range("C3").Value = 427502.9215053763
$ws.Range("D3").Value = 424099.6974731183
$ws.Range("E3").Value = 762.8306451612904
$ws.Range("F3").Value = 429900
$ws.Range("G3").Value = 430000
$ws.Range("H3").Value = 649
$ws.Range("I3").Value = 1860
$ws.Range("J3").Value = 555.96
$ws.Range("C4").Value = 426933.6651612903
$ws.Range("D4").Value = 422437.3082516129
$ws.Range("E4").Value = 789.2909677419354
$ws.Range("F4").Value = 429900
$ws.Range("G4").Value = 429000
$ws.Range("I4").Value = 1550
$ws.Range("J4").Value = 535.21

$ws = $wb.Worksheets.Item("2019")
$ws.Range("C2").Value = 428194.3145817913
$ws.Range("D2").Value = 425550.5336787565
$ws.Range("E2").Value = 792.4122871946706
$ws.Range("F2").Value = 438800
$ws.Range("G2").Value = 435000
$ws.Range("I2").Value = 1351
$ws.Range("J2").Value = 537.03
$ws.Range("C3").Value = 426447.5816266823
$ws.Range("D3").Value = 425909.9361497952
$ws.Range("E3").Value = 801.428320655354
$ws.Range("F3").Value = 432000
$ws.Range("G3").Value = 434000
$ws.Range("I3").Value = 1709
$ws.Range("J3").Value = 531.4400000000001
$ws.Range("C4").Value = 430891.138317757
$ws.Range("D4").Value = 430233.6279626168
$ws.Range("E4").Value = 798.5327102803739
$ws.Range("F4").Value = 439000
$ws.Range("G4").Value = 438588
$ws.Range("I4").Value = 1605
$ws.Range("J4").Value = 538.78
$ws.Range("C5").Value = 431979.1935483871
$ws.Range("D5").Value = 431607.3471986417
$ws.Range("E5").Value = 838.304753820034
$ws.Range("F5").Value = 439700
$ws.Range("G5").Value = 440000
$ws.Range("I5").Value = 1178
$ws.Range("J5").Value = 514.86

$ws = $wb.Worksheets.Item("2020")
$ws.Range("C2").Value = 426104.2031662269
$ws.Range("D2").Value = 438234.8443403694
$ws.Range("E2").Value = 832.8733509234828
$ws.Range("F2").Value = 429900
$ws.Range("G2").Value = 448650
$ws.Range("I2").Value = 758
$ws.Range("J2").Value = 526.17
$ws.Range("C3").Value = 440571.3618421053
$ws.Range("D3").Value = 439270.6101973684
$ws.Range("E3").Value = 776.2220394736842
$ws.Range("G3").Value = 448194
$ws.Range("I3").Value = 608
$ws.Range("J3").Value = 565.91
